# Automatische test-sync: 2025-06-17 23:04:10
# Appends the newly received mail-log entry as row 66 on the "Logs" sheet,
# extends the conditional formatting ranges to cover the new row, and
# updates the "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new log entry -------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 66
$logs.Range("A" + $newRow).Value = "Vragen over samenwerking"
$logs.Range("B" + $newRow).Value = "mailmind.test@zohomail.eu"
$logs.Range("C" + $newRow).Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D" + $newRow).Value = "Overig"
$logs.Range("F" + $newRow).Value = "2025-06-17 23:03:10"
$logs.Range("G" + $newRow).Value = "Nee"

# --- Extend the conditional formatting so it still covers the full table
$catRange = $logs.Range("D2:D65")
$catConditions = $catRange.FormatConditions
for ($i = 1; $i -le $catConditions.Count; $i++) {
    $catConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D66"))
}

$answeredRange = $logs.Range("G2:G65")
$answeredConditions = $answeredRange.FormatConditions
for ($i = 1; $i -le $answeredConditions.Count; $i++) {
    $answeredConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G66"))
}

# --- Dashboard sheet: bump the "Overig" counter --------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 19
